$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing wrap-text style from the data rows down onto the new
# blank rows (13-21) before we touch the data range itself.
$ws.Range("A2").Copy()
$ws.Range("A13:A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear existing data rows (A2:A11) formatting + content, we'll rewrite them
$ws.Range("A2:A11").Clear()

# New sorted list of work types including the newly added "202 Goods receive G"
$values = @(
    "202 Goods receive G",
    "203 Goods receive D",
    "206 Replenishment",
    "208 Sorter 2 SALE",
    "209 Sortation",
    "210 Pick Sortation",
    "211 Misc",
    "213 Pick3PL_Astro",
    "214 Large orders",
    "215 Picking Online",
    "217 Returns"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Select cell F21 to match the final selection state
$ws.Range("F21").Select()
